# Weekly update: insert a new week's observation for
# "Feria Lagunitas de Puerto Montt - Zanahoria" at the top of the data
# (row 272), pushing every subsequent row down by one and extending the
# series through row 325.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 272..324 down to 273..325, carrying all existing formatting
# (including the date-style column D) along with them.
$ws.Rows.Item(272).Insert()

# Populate the freshly inserted row 272 with the new weekly record.
$ws.Cells.Item(272, 1).Value  = 4
$ws.Cells.Item(272, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(272, 3).Value  = "Los Lagos"
$ws.Cells.Item(272, 4).Value  = 44637
$ws.Cells.Item(272, 5).Value  = 10
$ws.Cells.Item(272, 6).Value  = 100114013
$ws.Cells.Item(272, 7).Value  = "Zanahoria"
$ws.Cells.Item(272, 8).Value  = "Sin especificar"
$ws.Cells.Item(272, 9).Value  = "Primera"
$ws.Cells.Item(272, 10).Value = 250
$ws.Cells.Item(272, 11).Value = 11000
$ws.Cells.Item(272, 12).Value = 11000
$ws.Cells.Item(272, 13).Value = 11000
$ws.Cells.Item(272, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(272, 15).Value = "Chillán"
$ws.Cells.Item(272, 16).Value = 550
$ws.Cells.Item(272, 17).Value = 20
$ws.Cells.Item(272, 18).Value = "Hortaliza"
